$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save off the current (pre-shift) content of rows 148-149 so we can
# re-insert it, unchanged, two rows further down (new rows 150-151),
# before we overwrite 148-149 with their updated values.
$row148 = @($ws.Cells.Item(148, 1).Value2, $ws.Cells.Item(148, 2).Value2, $ws.Cells.Item(148, 3).Value2, $ws.Cells.Item(148, 4).Value2, $ws.Cells.Item(148, 5).Value2, $ws.Cells.Item(148, 6).Value2, $ws.Cells.Item(148, 7).Value2, $ws.Cells.Item(148, 8).Value2, $ws.Cells.Item(148, 9).Value2, $ws.Cells.Item(148, 10).Value2, $ws.Cells.Item(148, 11).Value2, $ws.Cells.Item(148, 12).Value2, $ws.Cells.Item(148, 13).Value2, $ws.Cells.Item(148, 14).Value2, $ws.Cells.Item(148, 15).Value2, $ws.Cells.Item(148, 16).Value2, $ws.Cells.Item(148, 17).Value2, $ws.Cells.Item(148, 18).Value2)
$row149 = @($ws.Cells.Item(149, 1).Value2, $ws.Cells.Item(149, 2).Value2, $ws.Cells.Item(149, 3).Value2, $ws.Cells.Item(149, 4).Value2, $ws.Cells.Item(149, 5).Value2, $ws.Cells.Item(149, 6).Value2, $ws.Cells.Item(149, 7).Value2, $ws.Cells.Item(149, 8).Value2, $ws.Cells.Item(149, 9).Value2, $ws.Cells.Item(149, 10).Value2, $ws.Cells.Item(149, 11).Value2, $ws.Cells.Item(149, 12).Value2, $ws.Cells.Item(149, 13).Value2, $ws.Cells.Item(149, 14).Value2, $ws.Cells.Item(149, 15).Value2, $ws.Cells.Item(149, 16).Value2, $ws.Cells.Item(149, 17).Value2, $ws.Cells.Item(149, 18).Value2)

# Insert two new rows right after row 149 (i.e. before old row 150),
# pushing the old rows 150-175 down to 152-177.
$ws.Rows.Item(150).Insert()
$ws.Rows.Item(150).Insert()

# Re-write the saved rows 148-149 content into the newly opened rows 150-151.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(150, $col).Value2 = $row148[$col - 1]
    $ws.Cells.Item(151, $col).Value2 = $row149[$col - 1]
}

# Apply the date-format style used throughout column D to the two new rows.
$ws.Range("D150").NumberFormat = $ws.Range("D152").NumberFormat
$ws.Range("D151").NumberFormat = $ws.Range("D152").NumberFormat

# Update rows 148-149 in place with the new survey data (later date, higher prices).
$ws.Cells.Item(148, 4).Value2 = 44855
$ws.Cells.Item(148, 11).Value2 = 700
$ws.Cells.Item(148, 12).Value2 = 800
$ws.Cells.Item(148, 13).Value2 = 750
$ws.Cells.Item(148, 16).Value2 = 750

$ws.Cells.Item(149, 4).Value2 = 44855
$ws.Cells.Item(149, 11).Value2 = 600
$ws.Cells.Item(149, 12).Value2 = 600
$ws.Cells.Item(149, 13).Value2 = 600
$ws.Cells.Item(149, 16).Value2 = 600
